# Update transition-probability matrix on "Washington St._B" sheet.
# Commit: "added more games, sped up simulate game logic, and drafted optimization logic"
# This updates the probabilities in rows 2-19 (columns B-S) that shifted as a
# result of additional simulated games being folded into the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1481481481481481
$ws.Cells.Item(2, 3).Value = 0.6407407407407407
$ws.Cells.Item(2, 10).Value = 0.01481481481481482
$ws.Cells.Item(2, 16).Value = 0.1185185185185185
$ws.Cells.Item(2, 19).Value = 0.07777777777777778
$ws.Cells.Item(3, 3).Value = 0.04395604395604396
$ws.Cells.Item(3, 10).Value = 0.01098901098901099
$ws.Cells.Item(3, 16).Value = 0.7032967032967034
$ws.Cells.Item(3, 19).Value = 0.2417582417582418
$ws.Cells.Item(4, 10).Value = 0.05263157894736842
$ws.Cells.Item(4, 16).Value = 0.7105263157894737
$ws.Cells.Item(4, 19).Value = 0.2368421052631579
$ws.Cells.Item(5, 16).Value = 0.5
$ws.Cells.Item(5, 19).Value = 0.5
$ws.Cells.Item(6, 2).Value = 0.0776255707762557
$ws.Cells.Item(6, 4).Value = 0.0182648401826484
$ws.Cells.Item(6, 6).Value = 0.0821917808219178
$ws.Cells.Item(6, 10).Value = 0.2328767123287671
$ws.Cells.Item(6, 15).Value = 0.0045662100456621
$ws.Cells.Item(6, 17).Value = 0.1506849315068493
$ws.Cells.Item(6, 18).Value = 0.0639269406392694
$ws.Cells.Item(6, 19).Value = 0.3698630136986301
$ws.Cells.Item(7, 2).Value = 0.1555555555555556
$ws.Cells.Item(7, 4).Value = 0.02777777777777778
$ws.Cells.Item(7, 6).Value = 0.05555555555555555
$ws.Cells.Item(7, 10).Value = 0.1333333333333333
$ws.Cells.Item(7, 15).Value = 0.005555555555555556
$ws.Cells.Item(7, 17).Value = 0.1888888888888889
$ws.Cells.Item(7, 18).Value = 0.07777777777777778
$ws.Cells.Item(7, 19).Value = 0.3555555555555556
$ws.Cells.Item(8, 2).Value = 0.08659793814432989
$ws.Cells.Item(8, 4).Value = 0.01443298969072165
$ws.Cells.Item(8, 6).Value = 0.05567010309278351
$ws.Cells.Item(8, 10).Value = 0.1092783505154639
$ws.Cells.Item(8, 15).Value = 0.02268041237113402
$ws.Cells.Item(8, 17).Value = 0.1938144329896907
$ws.Cells.Item(8, 18).Value = 0.06804123711340206
$ws.Cells.Item(8, 19).Value = 0.4494845360824742
$ws.Cells.Item(9, 2).Value = 0.06626506024096386
$ws.Cells.Item(9, 4).Value = 0.01807228915662651
$ws.Cells.Item(9, 6).Value = 0.07228915662650602
$ws.Cells.Item(9, 10).Value = 0.1144578313253012
$ws.Cells.Item(9, 15).Value = 0.03012048192771084
$ws.Cells.Item(9, 17).Value = 0.1686746987951807
$ws.Cells.Item(9, 18).Value = 0.1024096385542169
$ws.Cells.Item(9, 19).Value = 0.427710843373494
$ws.Cells.Item(10, 2).Value = 0.1216834400731931
$ws.Cells.Item(10, 4).Value = 0.01829826166514181
$ws.Cells.Item(10, 5).Value = 0.001829826166514181
$ws.Cells.Item(10, 6).Value = 0.08417200365965233
$ws.Cells.Item(10, 10).Value = 0.1033851784080512
$ws.Cells.Item(10, 15).Value = 0.01463860933211345
$ws.Cells.Item(10, 17).Value = 0.1957913998170174
$ws.Cells.Item(10, 18).Value = 0.0686184812442818
$ws.Cells.Item(10, 19).Value = 0.3915827996340348
$ws.Cells.Item(11, 7).Value = 0.1232394366197183
$ws.Cells.Item(11, 10).Value = 0.1373239436619718
$ws.Cells.Item(11, 11).Value = 0.176056338028169
$ws.Cells.Item(11, 12).Value = 0.5492957746478874
$ws.Cells.Item(11, 19).Value = 0.01408450704225352
$ws.Cells.Item(12, 7).Value = 0.7453416149068323
$ws.Cells.Item(12, 10).Value = 0.1801242236024845
$ws.Cells.Item(12, 12).Value = 0.02484472049689441
$ws.Cells.Item(12, 19).Value = 0.04968944099378882
$ws.Cells.Item(13, 7).Value = 0.7021276595744681
$ws.Cells.Item(13, 10).Value = 0.2978723404255319
$ws.Cells.Item(14, 7).Value = 0.6666666666666666
$ws.Cells.Item(14, 10).Value = 0.3333333333333333
$ws.Cells.Item(15, 6).Value = 0.004830917874396135
$ws.Cells.Item(15, 8).Value = 0.178743961352657
$ws.Cells.Item(15, 9).Value = 0.08695652173913043
$ws.Cells.Item(15, 10).Value = 0.3140096618357488
$ws.Cells.Item(15, 11).Value = 0.1159420289855072
$ws.Cells.Item(15, 15).Value = 0.05314009661835749
$ws.Cells.Item(15, 19).Value = 0.2463768115942029
$ws.Cells.Item(16, 6).Value = 0.03314917127071823
$ws.Cells.Item(16, 8).Value = 0.1823204419889503
$ws.Cells.Item(16, 9).Value = 0.06629834254143646
$ws.Cells.Item(16, 10).Value = 0.3701657458563536
$ws.Cells.Item(16, 11).Value = 0.1049723756906077
$ws.Cells.Item(16, 13).Value = 0.01104972375690608
$ws.Cells.Item(16, 15).Value = 0.06629834254143646
$ws.Cells.Item(16, 19).Value = 0.1657458563535912
$ws.Cells.Item(17, 6).Value = 0.02
$ws.Cells.Item(17, 8).Value = 0.225
$ws.Cells.Item(17, 9).Value = 0.08500000000000001
$ws.Cells.Item(17, 10).Value = 0.385
$ws.Cells.Item(17, 11).Value = 0.08500000000000001
$ws.Cells.Item(17, 13).Value = 0.0175
$ws.Cells.Item(17, 14).Value = 0.0025
$ws.Cells.Item(17, 15).Value = 0.055
$ws.Cells.Item(17, 19).Value = 0.125
$ws.Cells.Item(18, 6).Value = 0.03947368421052631
$ws.Cells.Item(18, 8).Value = 0.2302631578947368
$ws.Cells.Item(18, 9).Value = 0.09868421052631579
$ws.Cells.Item(18, 10).Value = 0.3486842105263158
$ws.Cells.Item(18, 11).Value = 0.09868421052631579
$ws.Cells.Item(18, 13).Value = 0.03289473684210526
$ws.Cells.Item(18, 14).Value = 0.006578947368421052
$ws.Cells.Item(18, 15).Value = 0.04605263157894737
$ws.Cells.Item(18, 19).Value = 0.09868421052631579
$ws.Cells.Item(19, 6).Value = 0.01130856219709208
$ws.Cells.Item(19, 8).Value = 0.2374798061389338
$ws.Cells.Item(19, 9).Value = 0.07189014539579967
$ws.Cells.Item(19, 10).Value = 0.3392568659127625
$ws.Cells.Item(19, 11).Value = 0.1122778675282714
$ws.Cells.Item(19, 13).Value = 0.02584814216478191
$ws.Cells.Item(19, 14).Value = 0.0008077544426494346
$ws.Cells.Item(19, 15).Value = 0.07996768982229402
$ws.Cells.Item(19, 19).Value = 0.1211631663974152
